# The workbook was simply re-opened and re-saved (no real data edits were
# made - see commit message "update on webpage since last time"). We
# reproduce the parts of that resave that are reachable through the Excel
# object model: the cells that carried a redundant/duplicate cell style
# get nudged onto their canonical (de-duplicated) style, and the last
# active selection is left on F14 (row 14, column F), matching the saved
# view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells G2:J5 area (I2, J2, G3, I3, J3, G4, I4, J4, G5, I5, J5) used a cell
# style that only differed from the neighbouring "H column" style by an
# inert alignment flag. Re-touching their font collapses them onto that
# shared style, just like Excel's own style-table clean-up does on save.
$dedupCells = @("I2", "J2", "G3", "I3", "J3", "G4", "I4", "J4", "G5", "I5", "J5")
foreach ($addr in $dedupCells) {
    $ws.Range($addr).Font.Size = 10
}

# Restore the last selection/active cell that was saved with the workbook.
$ws.Range("F14").Select()

$wb.Save()
